$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new data row at row 219 (pushes existing rows 219.. down by one) ---
$ws.Rows.Item(219).Insert()
$ws.Cells.Item(219,1).Value  = 7
$ws.Cells.Item(219,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(219,3).Value  = "Ñuble"
$ws.Cells.Item(219,4).Value  = 45006
$ws.Cells.Item(219,5).Value  = 16
$ws.Cells.Item(219,6).Value  = 100112043
$ws.Cells.Item(219,7).Value  = "Pepino ensalada"
$ws.Cells.Item(219,8).Value  = "Sin especificar"
$ws.Cells.Item(219,9).Value  = "Primera"
$ws.Cells.Item(219,10).Value = 80
$ws.Cells.Item(219,11).Value = 15000
$ws.Cells.Item(219,12).Value = 15000
$ws.Cells.Item(219,13).Value = 15000
$ws.Cells.Item(219,14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(219,15).Value = "Región del Maule"
$ws.Cells.Item(219,16).Value = 188
$ws.Cells.Item(219,17).Value = 80
$ws.Cells.Item(219,18).Value = "Hortaliza"

# --- Insert a second new data row at (the now-shifted) row 259 ---
$ws.Rows.Item(259).Insert()
$ws.Cells.Item(259,1).Value  = 7
$ws.Cells.Item(259,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(259,3).Value  = "Ñuble"
$ws.Cells.Item(259,4).Value  = 45005
$ws.Cells.Item(259,5).Value  = 16
$ws.Cells.Item(259,6).Value  = 100112043
$ws.Cells.Item(259,7).Value  = "Pepino ensalada"
$ws.Cells.Item(259,8).Value  = "Sin especificar"
$ws.Cells.Item(259,9).Value  = "Primera"
$ws.Cells.Item(259,10).Value = 60
$ws.Cells.Item(259,11).Value = 15000
$ws.Cells.Item(259,12).Value = 15000
$ws.Cells.Item(259,13).Value = 15000
$ws.Cells.Item(259,14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(259,15).Value = "Región del Maule"
$ws.Cells.Item(259,16).Value = 188
$ws.Cells.Item(259,17).Value = 80
$ws.Cells.Item(259,18).Value = "Hortaliza"
